$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1835
$ws.Range("F3").Value = 19
$ws.Range("F4").Value = 24
$ws.Range("F5").Value = 41
$ws.Range("F8").Value = 179
$ws.Range("F9").Value = 591
$ws.Range("F10").Value = 61
$ws.Range("F11").Value = 455
$ws.Range("F12").Value = 535
$ws.Range("F13").Value = 1403
$ws.Range("F14").Value = 1208
$ws.Range("F15").Value = 1420
$ws.Range("F16").Value = 19
$ws.Range("F17").Value = 1169
$ws.Range("F18").Value = 282
$ws.Range("F19").Value = 1564
$ws.Range("F20").Value = 439
$ws.Range("F21").Value = 769
$ws.Range("F25").Value = 1271
$ws.Range("F27").Value = 43
$ws.Range("F28").Value = 791
$ws.Range("F29").Value = 543
$ws.Range("F30").Value = 1026
$ws.Range("F31").Value = 238490
$ws.Range("F32").Value = 949
$ws.Range("F33").Value = 548
$ws.Range("F35").Value = 881
$ws.Range("F36").Value = 1037
$ws.Range("F37").Value = 23
$ws.Range("F38").Value = 823
$ws.Range("F39").Value = 1561
$ws.Range("F40").Value = 87
$ws.Range("F41").Value = 26
$ws.Range("F42").Value = 779
$ws.Range("F44").Value = 765
$ws.Range("F45").Value = 106

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 120
$ws.Range("F6").Value = 139
$ws.Range("F10").Value = 157
$ws.Range("F11").Value = 1401
$ws.Range("F12").Value = 65
$ws.Range("F14").Value = 2492
$ws.Range("F15").Value = 1173
$ws.Range("F16").Value = 386
$ws.Range("F19").Value = 29
$ws.Range("F20").Value = 63
$ws.Range("F23").Value = 421
$ws.Range("F24").Value = 20
$ws.Range("F26").Value = 271
$ws.Range("F27").Value = 58859
$ws.Range("F34").Value = 54
$ws.Range("F37").Value = 52
$ws.Range("F43").Value = 35
$ws.Range("F44").Value = 35
$ws.Range("F46").Value = 116
$ws.Range("F47").Value = 52

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 2681
$ws.Range("F7").Value = 4486
$ws.Range("F8").Value = 115
$ws.Range("F10").Value = 520
$ws.Range("F11").Value = 616
$ws.Range("F12").Value = 407
$ws.Range("F13").Value = 172
$ws.Range("F14").Value = 650
$ws.Range("F15").Value = 182
$ws.Range("F16").Value = 334

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1835
$ws.Range("F4").Value = 19
$ws.Range("F5").Value = 4486
$ws.Range("F6").Value = 616
$ws.Range("F7").Value = 41
$ws.Range("F8").Value = 172
$ws.Range("F9").Value = 172
$ws.Range("F10").Value = 650
$ws.Range("F11").Value = 650
$ws.Range("F12").Value = 182
$ws.Range("F13").Value = 139
$ws.Range("F16").Value = 179
$ws.Range("F17").Value = 157
$ws.Range("F18").Value = 1401
$ws.Range("F19").Value = 591
$ws.Range("F20").Value = 455
$ws.Range("F21").Value = 535
$ws.Range("F22").Value = 2492
$ws.Range("F23").Value = 1173
$ws.Range("F24").Value = 1403
$ws.Range("F25").Value = 1208
$ws.Range("F26").Value = 1420
$ws.Range("F27").Value = 1169
$ws.Range("F29").Value = 63
$ws.Range("F30").Value = 1564
$ws.Range("F31").Value = 769
$ws.Range("F33").Value = 334
$ws.Range("F34").Value = 421
$ws.Range("F35").Value = 1271
$ws.Range("F36").Value = 791
$ws.Range("F37").Value = 543
$ws.Range("F38").Value = 1026
$ws.Range("F39").Value = 271
$ws.Range("F40").Value = 949
$ws.Range("F41").Value = 881
$ws.Range("F42").Value = 23
$ws.Range("F43").Value = 823
$ws.Range("F45").Value = 1561
$ws.Range("F46").Value = 87
$ws.Range("F49").Value = 779
$ws.Range("F50").Value = 35
$ws.Range("F51").Value = 765
